# Junction_Flooding_444.xlsx — refresh sampled data (new timestamps/values for
# rows 2-5), drop the now-unused last data row (row 6), and widen the data
# columns (B:AH) to fit the new values ("custom accuracy + 데이터 1000개").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New data for rows 2-5 (column A = timestamp, B..AH = junction values) ---
$row2 = @(45125.50694444445,5.885,1.942,1.363,7.767,3.737,0.96,6.281,2.168,0.758,1.218,2.359,5.493,0.668,0.884,2.495,1.484,1.256,0,25.797,5.514,3.303,3.938,2.803,0.246,1.46,1.634,0.588,2.048,3.034,0.186,2.547,0.594,1.826)
$row3 = @(45125.51388888889,16.094,11.262,1.015,33.1,26.3,11.068,39.373,17.483,7.667,11.54,12.984,14.91,3.735,10.887,16.511,9.534000000000001,0.718,0.298,169.635,32.347,11.195,22.155,11.889,1.519,19.551,9.467000000000001,7.783,9.933999999999999,13.9,0.173,34.298,5.828,13.092)
$row4 = @(45125.52083333334,1.71,0.8169999999999999,0.346,2.549,1.444,0.328,9.874000000000001,0.827,0.403,0.398,0.834,1.588,0.282,0.21,1.264,0.5629999999999999,0.366,0,4.034,2.647,0.958,2.214,1.145,0.095,4.103,0.599,0.183,0.737,1.038,0.141,8.864000000000001,0.187,0.654)
$row5 = @(45125.52777777778,1.5,0.8100000000000001,0.25,2.43,1.53,0.43,5.48,0.93,0.29,0.49,0.84,1.4,0.24,0.33,1.07,0.59,0.27,0,3.93,2.18,0.88,1.62,0.92,0.1,2.07,0.61,0.26,0.72,1,0.12,4.51,0.24,0.7)

$dataRows = @(2,3,4,5)
$rowsData = @($row2,$row3,$row4,$row5)

for ($ri = 0; $ri -lt $dataRows.Length; $ri++) {
    $r = $dataRows[$ri]
    $vals = $rowsData[$ri]
    for ($ci = 0; $ci -lt $vals.Length; $ci++) {
        $col = $ci + 1
        $ws.Cells.Item($r, $col).Value = $vals[$ci]
    }
}

# --- 2. Drop row 6 (data now ends at row 5) ---
$ws.Rows.Item(6).Delete()

# --- 3. Resize columns B..AH (2..34) to the new widths ---
$newWidths = @(8,8,7,7,7,8,8,8,7,7,8,7,7,8,8,7,7,7,9,8,8,8,8,7,8,7,7,7,7,7,8,7,8)
for ($i = 0; $i -lt $newWidths.Length; $i++) {
    $col = $i + 2
    $ws.Columns.Item($col).ColumnWidth = $newWidths[$i] - 0.8333333333333333
}
